$d = $word.ActiveDocument

$d.Content.Find.Execute("22×18=396", $true, $false, $false, $false, $false, $true, 1, $false, "31×70=2170", 2) | Out-Null
$d.Content.Find.Execute("35×64=2240", $true, $false, $false, $false, $false, $true, 1, $false, "53×89=4717", 2) | Out-Null
$d.Content.Find.Execute("31×33=1023", $true, $false, $false, $false, $false, $true, 1, $false, "26×60=1560", 2) | Out-Null
$d.Content.Find.Execute("39×56=2184", $true, $false, $false, $false, $false, $true, 1, $false, "24×89=2136", 2) | Out-Null
$d.Content.Find.Execute("23×31=713", $true, $false, $false, $false, $false, $true, 1, $false, "33×95=3135", 2) | Out-Null
$d.Content.Find.Execute("30×74=2220", $true, $false, $false, $false, $false, $true, 1, $false, "98×12=1176", 2) | Out-Null
$d.Content.Find.Execute("70×70=4900", $true, $false, $false, $false, $false, $true, 1, $false, "61×30=1830", 2) | Out-Null
$d.Content.Find.Execute("41×54=2214", $true, $false, $false, $false, $false, $true, 1, $false, "75×87=6525", 2) | Out-Null
$d.Content.Find.Execute("64×39=2496", $true, $false, $false, $false, $false, $true, 1, $false, "18×91=1638", 2) | Out-Null
$d.Content.Find.Execute("85×27=2295", $true, $false, $false, $false, $false, $true, 1, $false, "74×38=2812", 2) | Out-Null
$d.Content.Find.Execute("73×60=4380", $true, $false, $false, $false, $false, $true, 1, $false, "53×93=4929", 2) | Out-Null
$d.Content.Find.Execute("95×33=3135", $true, $false, $false, $false, $false, $true, 1, $false, "35×99=3465", 2) | Out-Null
$d.Content.Find.Execute("88×16=1408", $true, $false, $false, $false, $false, $true, 1, $false, "72×64=4608", 2) | Out-Null
$d.Content.Find.Execute("20×79=1580", $true, $false, $false, $false, $false, $true, 1, $false, "48×11=528", 2) | Out-Null
$d.Content.Find.Execute("20×88=1760", $true, $false, $false, $false, $false, $true, 1, $false, "34×11=374", 2) | Out-Null
$d.Content.Find.Execute("73×84=6132", $true, $false, $false, $false, $false, $true, 1, $false, "59×88=5192", 2) | Out-Null
$d.Content.Find.Execute("96×79=7584", $true, $false, $false, $false, $false, $true, 1, $false, "19×40=760", 2) | Out-Null
$d.Content.Find.Execute("96×76=7296", $true, $false, $false, $false, $false, $true, 1, $false, "19×87=1653", 2) | Out-Null
$d.Content.Find.Execute("30×13=390", $true, $false, $false, $false, $false, $true, 1, $false, "29×57=1653", 2) | Out-Null
$d.Content.Find.Execute("89×11=979", $true, $false, $false, $false, $false, $true, 1, $false, "59×70=4130", 2) | Out-Null
$d.Content.Find.Execute("56×58=3248", $true, $false, $false, $false, $false, $true, 1, $false, "28×97=2716", 2) | Out-Null
$d.Content.Find.Execute("56×87=4872", $true, $false, $false, $false, $false, $true, 1, $false, "97×23=2231", 2) | Out-Null
$d.Content.Find.Execute("68×66=4488", $true, $false, $false, $false, $false, $true, 1, $false, "72×37=2664", 2) | Out-Null
$d.Content.Find.Execute("61×78=4758", $true, $false, $false, $false, $false, $true, 1, $false, "87×23=2001", 2) | Out-Null
$d.Content.Find.Execute("76×69=5244", $true, $false, $false, $false, $false, $true, 1, $false, "26×14=364", 2) | Out-Null
